# Add a new patient record (row 6) to the Patients sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First name / last name are plain text - straightforward assignment.
$ws.Range("A6").Value = "Piotr"
$ws.Range("B6").Value = "Kowalski"

# The PESEL number has a leading zero ("0430403") and must be stored as
# text, otherwise Excel would coerce it to a number and the leading zero
# would be lost. Forcing a text number format achieves that, but it also
# allocates a cell style for C6 that the source workbook never had.
# Re-applying a style-free format (copied from an unstyled cell) via
# PasteSpecial(values) afterwards keeps the text content while dropping
# the extra style, so the cell ends up as a plain shared-string cell,
# matching how the sheet was originally authored.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "0430403"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
